$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Q0) ---
$ws.Range("B2").Value = -0.01636740676293108
$ws.Range("C2").Value = 0.7208754989987691
$ws.Range("D2").Value = 0.8155649038874973
$ws.Range("E2").Value = 0.9030863213931973
$ws.Range("F2").Value = 0.9370230342118662
$ws.Range("G2").Value = 14

# --- Row 3 (Q1) ---
$ws.Range("B3").Value = 0.1136094414023552
$ws.Range("C3").Value = 0.6727610774151302
$ws.Range("D3").Value = 0.6903787083050854
$ws.Range("E3").Value = 0.8308903106337716
$ws.Range("F3").Value = 0.8566957278151761
$ws.Range("G3").Value = 13

# --- Row 4 (Q2) ---
$ws.Range("B4").Value = 0.1783636749128446
$ws.Range("C4").Value = 0.6586562652618152
$ws.Range("D4").Value = 0.5234084977275197
$ws.Range("E4").Value = 0.7234697628287721
$ws.Range("F4").Value = 0.7323150567884177
$ws.Range("G4").Value = 12

# --- Row 5 (Q3) ---
$ws.Range("B5").Value = 0.240625490472168
$ws.Range("C5").Value = 0.716344199269209
$ws.Range("D5").Value = 0.6084117250097457
$ws.Range("E5").Value = 0.7800075159956766
$ws.Range("F5").Value = 0.7781787764898576
$ws.Range("G5").Value = 11

# --- Row 6 (Q4) ---
$ws.Range("B6").Value = 0.2137183226789437
$ws.Range("C6").Value = 0.6340062147284611
$ws.Range("D6").Value = 0.4822391445073304
$ws.Range("E6").Value = 0.6944344061949482
$ws.Range("F6").Value = 0.6964701661143612
$ws.Range("G6").Value = 10

# --- Row 7 (Q5) ---
$ws.Range("B7").Value = 0.1237973885243304
$ws.Range("C7").Value = 0.583296279823095
$ws.Range("D7").Value = 0.3672486758207861
$ws.Range("E7").Value = 0.6060104585077605
$ws.Range("F7").Value = 0.6292163719399392
$ws.Range("G7").Value = 9

# --- Row 8 (Q6) ---
$ws.Range("B8").Value = 0.3031634011949809
$ws.Range("C8").Value = 0.5793847622031741
$ws.Range("D8").Value = 0.4383022118985093
$ws.Range("E8").Value = 0.6620439652307913
$ws.Range("F8").Value = 0.6447270716274294
$ws.Range("G8").Value = 6

# --- Row 9 (Q7) ---
$ws.Range("B9").Value = 0.08312183816600854
$ws.Range("C9").Value = 0.1419465487444733
$ws.Range("D9").Value = 0.02590123319020616
$ws.Range("E9").Value = 0.1609386006842552
$ws.Range("F9").Value = 0.1687838553155042
$ws.Range("G9").Value = 3

# --- Row 10 (new row, Q8) ---
# Copy formatting from the A9 label cell (bold/centered/bordered style)
# so the new A10 label matches the existing Q-row labels.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Q8"

$ws.Range("B10").Value = -0.7701202119308102
$ws.Range("C10").Value = 0.7701202119308102
$ws.Range("D10").Value = 0.5930851408243559
$ws.Range("E10").Value = 0.7701202119308102
$ws.Range("G10").Value = 1
